$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.367.24'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '2.103.61'
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '229.18'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.35%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.614'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.64%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '61.28'
$c.ClearFormats()
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.48%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0847'
$c.ClearFormats()
$ws.Range("E10").Value = '  +3.32%  '
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").Value = '2.413.99'
$ws.Range("E12").Value = '  +3.15%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.79'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.94%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '22.41'
$c.ClearFormats()
$ws.Range("E14").Value = '  +5.94%  '
$ws.Range("E15").Value = '  +5.93%  '
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = '2.063.77'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '38.281.58'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("E19").Value = '  +2.12%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '70.48'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").Value = '  +1.33%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '224.66'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.80%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.ClearFormats()
$ws.Range("E25").Value = '  +2.68%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '169.87'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.22%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.44'
$c.ClearFormats()
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("E30").Value = '  +5.84%  '
$ws.Range("E31").Value = '  -0.66%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.ClearFormats()
$ws.Range("E32").Value = '  +8.22%  '
$ws.Range("E33").Value = '  +4.27%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.45'
$c.ClearFormats()
$ws.Range("E34").Value = '  +0.76%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0607'
$c.ClearFormats()
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.44'
$c.ClearFormats()
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.ClearFormats()
$ws.Range("E37").Value = '  +4.62%  '
$ws.Range("E38").Value = '  +5.71%  '
$ws.Range("E39").Value = '  +0.03%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.12'
$c.ClearFormats()
$ws.Range("E40").Value = '  +2.08%  '
$ws.Range("D41").Value = '1.550.79'
$ws.Range("E41").Value = '  +0.72%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '100.18'
$c.ClearFormats()
$ws.Range("E42").Value = '  +3.72%  '
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("E46").Value = '  +4.10%  '
$ws.Range("E47").Value = '  +1.38%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.51'
$c.ClearFormats()
$ws.Range("E48").Value = '  +4.65%  '
$ws.Range("E49").Value = '  +1.89%  '
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = '2.300.53'
$ws.Range("E51").Value = '  +3.20%  '
